$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 186 - this shifts existing rows 186..273 down to 187..274
$ws.Rows.Item(186).Insert()

# Populate the newly inserted row 186 with the new record
$ws.Cells.Item(186, 1).Value = 4
$ws.Cells.Item(186, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(186, 3).Value = "Los Lagos"
$ws.Cells.Item(186, 4).Value = 44609
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 100112008
$ws.Cells.Item(186, 7).Value = "Coliflor"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Segunda"
$ws.Cells.Item(186, 10).Value = 250
$ws.Cells.Item(186, 11).Value = 1500
$ws.Cells.Item(186, 12).Value = 1500
$ws.Cells.Item(186, 13).Value = 1500
$ws.Cells.Item(186, 14).Value = "`$/unidad"
$ws.Cells.Item(186, 15).Value = "Región Metropolitana"
$ws.Cells.Item(186, 16).Value = 1500
$ws.Cells.Item(186, 17).Value = 1
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (style index carried over from the row above)
$ws.Cells.Item(186, 4).NumberFormat = $ws.Cells.Item(185, 4).NumberFormat
